{"js": "// Docx writer: include abstract title (localized term for \"Abstract\").\n//\n// 1. Insert a new paragraph styled \"AbstractTitle\" containing the text\n//    \"Abstract\" right after the Author paragraph (\"A. M.\") and right\n//    before the existing \"Abstract\" (body) paragraph.\n// 2. Give the syntax-highlighting token styles `ImportTok` and\n//    `BuiltInTok` their green color (and, for ImportTok, bold) so import\n//    statements / built-ins render consistently with the rest of the\n//    Pandoc-generated token palette.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\n// Locate the \"Author\" paragraph (\"A. M.\") so the new title is inserted\n// immediately after it (and therefore immediately before the existing\n// \"Abstract\" paragraph).\nlet authorParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.style === \"Author\") {\n    authorParagraph = p;\n    break;\n  }\n}\n\nif (authorParagraph) {\n  const abstractTitle = authorParagraph.insertParagraph(\"Abstract\", \"After\");\n  abstractTitle.style = \"AbstractTitle\";\n}\n\n// Character styles used to highlight source-code tokens: add the shared\n// green token color, and (for ImportTok only) bold.\nconst styles = context.document.getStyles();\n\nconst importTok = styles.getByNameOrNullObject(\"ImportTok\");\nconst builtInTok = styles.getByNameOrNullObject(\"BuiltInTok\");\nawait context.sync();\n\nif (!importTok.isNullObject) {\n  importTok.font.color = \"#008000\";\n  importTok.font.bold = true;\n}\n\nif (!builtInTok.isNullObject) {\n  builtInTok.font.color = \"#008000\";\n}\n\nawait context.sync();\n", "ps1": "# Docx writer: include abstract title (localized term for \"Abstract\").\n#\n# 1. Insert a new paragraph styled \"AbstractTitle\" containing the text\n#    \"Abstract\" right after the Author paragraph (\"A. M.\") and right\n#    before the existing \"Abstract\" (body) paragraph.\n# 2. Give the syntax-highlighting token styles `ImportTok` and\n#    `BuiltInTok` their green color (and, for ImportTok, bold) so import\n#    statements / built-ins render consistently with the rest of the\n#    Pandoc-generated token palette.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Author\" paragraph (\"A. M.\") so the new title is inserted\n# immediately after it (and therefore immediately before the existing\n# \"Abstract\" paragraph).\n$authorParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Author\") {\n        $authorParagraph = $p\n        break\n    }\n}\n\nif ($authorParagraph -ne $null) {\n    $r = $authorParagraph.Range\n    $r.Collapse(0)  # wdCollapseEnd\n    $r.InsertParagraphAfter()\n\n    $newPara = $authorParagraph.Next()\n    $newPara.Range.Text = \"Abstract\"\n    $newPara.Style = \"AbstractTitle\"\n}\n\n# Character styles used to highlight source-code tokens: add the shared\n# green token color, and (for ImportTok only) bold.\n$importTok = $d.Styles(\"ImportTok\")\n$importTok.Font.Color = 32768      # wdColorGreen (RGB 0x008000)\n$importTok.Font.Bold = 1\n\n$builtInTok = $d.Styles(\"BuiltInTok\")\n$builtInTok.Font.Color = 32768     # wdColorGreen (RGB 0x008000)\n"}
